$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isyatirim")

# New rows of data to append (Tarih, Kapanis, Min, Max, AOF, Hacim, Sermaye, USDTRY, BIST100, PiyasaDegeriTL, PiyasaDegeriUSD, HalkaAcikTL, HalkaAcikUSD)
$rows = @(
    @("08-11-2023", 318.25, 314.5, 325, 318.89, 425193518, 124, 28.528, 7861, 39383, 1381, 9129, 320),
    @("09-11-2023", 316.5, 311, 322, 315.64, 499649766, 124, 28.5272, 7841, 39167, 1373, 9079, 318),
    @("10-11-2023", 311, 309, 318.5, 312.62, 318445923, 124, 28.537, 7771, 38486, 1349, 8921, 313)
)

$startRow = 120
$endRow = $startRow + $rows.Count - 1

# Keep the date column as text (matching the rest of the "Tarih" column)
# so Excel does not reinterpret strings like "08-11-2023" as date serials.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $data[$c]
    }
}
